$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.036.00'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.952.83'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.89'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.36'
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.950.45'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.09'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('E11').Value = '  +6.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.440'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('E13').Value = '  +4.69%  '
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.444.58'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.047.42'
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.953.68'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '441.10'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.48'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.668'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.01'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.25'
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.82'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.83'
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.31'
$ws.Range('E30').Value = '  +6.12%  '
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0000103'
$ws.Range('E32').Value = '  +16.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.44'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.992'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.10'
$ws.Range('E37').Value = '  +4.26%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.61'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.75'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.49'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -3.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.280'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.54'
$ws.Range('E44').Value = '  -6.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.709.11'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.85'
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '361.29'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.104'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.77'
$ws.Range('E51').Value = '  -3.24%  '
